$d = $word.ActiveDocument

# 1. Update the letter date.
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the "2698 Greenrock Road, Milpitas CA 95035" paragraph into two
#    paragraphs: "2698 Greenrock Road" and a new "Milpitas, CA 95035" line,
#    mirroring the formatting of the original paragraph.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "2698 Greenrock Road, Milpitas CA 95035`r") {
        $r = $d.Range($p.Range.Start, $p.Range.End - 1)
        $r.Text = "2698 Greenrock Road"

        $p = $d.Paragraphs.Item($i)
        $p.Range.InsertParagraphAfter()

        $newPara = $d.Paragraphs.Item($i + 1)
        $newRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
        $newRange.Text = "Milpitas, CA 95035"
        break
    }
}

# 3. Remove the blank "No Spacing" paragraph that follows
#    "...Board of Directors".
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Board of Directors") {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text -eq "`r") {
            $next.Range.Delete()
        }
        break
    }
}
